# Aufgabenverteilung.xlsx — remove the "Scrum Modell" task row from the
# Tabelle1 plan (added properties/testApp/documents row shifted up).
#
# The row "Scrum Modell" (row 5) is removed entirely; Excel's native
# Delete-Row behaviour takes care of shifting every row below it up by
# one (values, styles, number formats, etc. all move together), which
# is exactly what the target workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Delete the whole row 5 ("Scrum Modell") - everything below shifts up.
$ws.Rows("5:5").Delete()

# Excel leaves behind a selection on the row that used to be row 6
# (now row 5) after this kind of edit - mirror that in the saved view.
$ws.Rows("5:5").Select() | Out-Null
